# Updates odds values in Sheet1 to match the latest FlashScore scrape
# (rows 3,4,5,6,7,11,12,18,23 — betting-odds/score columns only).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF3").Value = 6
$ws.Range("AK3").Value = 15
$ws.Range("AP3").Value = 1.88
$ws.Range("AQ3").Value = 2.02
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.63
$ws.Range("AA4").Value = 9.5
$ws.Range("AB4").Value = 21
$ws.Range("AJ4").Value = 9
$ws.Range("AK4").Value = 17
$ws.Range("AL4").Value = 13
$ws.Range("G4").Value = 2.25
$ws.Range("I4").Value = 3.4
$ws.Range("J4").Value = 3
$ws.Range("L4").Value = 4
$ws.Range("Z4").Value = 10
$ws.Range("AP5").Value = 2.15
$ws.Range("AQ5").Value = 1.67
$ws.Range("O5").Value = 1.62
$ws.Range("P5").Value = 2.2
$ws.Range("S5").Value = 6.5
$ws.Range("T5").Value = 1.11
$ws.Range("AB6").Value = 19
$ws.Range("AJ6").Value = 8.5
$ws.Range("AP6").Value = 2.03
$ws.Range("AQ6").Value = 1.83
$ws.Range("AR6").Value = 4.3
$ws.Range("AS6").Value = 1.21
$ws.Range("G6").Value = 2.1
$ws.Range("I6").Value = 4.1
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 6
$ws.Range("O6").Value = 1.53
$ws.Range("P6").Value = 2.38
$ws.Range("Q6").Value = 2.7
$ws.Range("R6").Value = 1.44
$ws.Range("U6").Value = 1.62
$ws.Range("V6").Value = 2.2
$ws.Range("Z6").Value = 8.5
$ws.Range("AA7").Value = 9.5
$ws.Range("AC7").Value = 17
$ws.Range("AE7").Value = 7
$ws.Range("AG7").Value = 29
$ws.Range("AH7").Value = 126
$ws.Range("AJ7").Value = 11
$ws.Range("AP7").Value = 1.93
$ws.Range("AQ7").Value = 1.93
$ws.Range("AR7").Value = 3.85
$ws.Range("AS7").Value = 1.25
$ws.Range("H7").Value = 3.8
$ws.Range("I7").Value = 6.5
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 7.5
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 1.5
$ws.Range("P7").Value = 2.5
$ws.Range("Q7").Value = 2.5
$ws.Range("R7").Value = 1.5
$ws.Range("S7").Value = 5
$ws.Range("T7").Value = 1.17
$ws.Range("U7").Value = 1.53
$ws.Range("V7").Value = 2.38
$ws.Range("W7").Value = 2.63
$ws.Range("X7").Value = 1.44
$ws.Range("Y7").Value = 4.75
$ws.Range("Z7").Value = 5.5
$ws.Range("AM11").Value = 23
$ws.Range("AP11").Value = 1.85
$ws.Range("AQ11").Value = 2
$ws.Range("G11").Value = 3.1
$ws.Range("H11").Value = 3.1
$ws.Range("I11").Value = 2.38
$ws.Range("L11").Value = 3.2
$ws.Range("O11").Value = 1.5
$ws.Range("P11").Value = 2.63
$ws.Range("Q11").Value = 2.5
$ws.Range("R11").Value = 1.53
$ws.Range("AB12").Value = 34
$ws.Range("AE12").Value = 7.5
$ws.Range("AI12").Value = 301
$ws.Range("AL12").Value = 10
$ws.Range("AM12").Value = 23
$ws.Range("AN12").Value = 21
$ws.Range("AO12").Value = 34
$ws.Range("G12").Value = 3.1
$ws.Range("H12").Value = 2.8
$ws.Range("I12").Value = 2.3
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 3.2
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 7.5
$ws.Range("O12").Value = 1.36
$ws.Range("P12").Value = 3
$ws.Range("Q12").Value = 2.25
$ws.Range("R12").Value = 1.62
$ws.Range("S12").Value = 4
$ws.Range("T12").Value = 1.22
$ws.Range("U12").Value = 1.5
$ws.Range("V12").Value = 2.5
$ws.Range("W12").Value = 1.83
$ws.Range("X12").Value = 1.83
$ws.Range("Y12").Value = 9
$ws.Range("Z12").Value = 15
$ws.Range("AC18").Value = 29
$ws.Range("AE18").Value = 6.5
$ws.Range("AH18").Value = 67
$ws.Range("AI18").Value = 501
$ws.Range("AJ18").Value = 6.5
$ws.Range("AK18").Value = 10
$ws.Range("AO18").Value = 41
$ws.Range("AP18").Value = 1.88
$ws.Range("AQ18").Value = 1.98
$ws.Range("G18").Value = 3.2
$ws.Range("H18").Value = 2.9
$ws.Range("I18").Value = 2.4
$ws.Range("J18").Value = 4
$ws.Range("M18").Value = 1.11
$ws.Range("N18").Value = 6.5
$ws.Range("O18").Value = 1.44
$ws.Range("P18").Value = 2.63
$ws.Range("Q18").Value = 2.5
$ws.Range("R18").Value = 1.5
$ws.Range("S18").Value = 5
$ws.Range("T18").Value = 1.17
$ws.Range("U18").Value = 1.57
$ws.Range("V18").Value = 2.25
$ws.Range("W18").Value = 2.05
$ws.Range("X18").Value = 1.7
$ws.Range("Z18").Value = 15
$ws.Range("AJ23").Value = 8
$ws.Range("AK23").Value = 17
$ws.Range("AP23").Value = 2.1
$ws.Range("AQ23").Value = 1.78
$ws.Range("G23").Value = 2.1
$ws.Range("I23").Value = 3.9
$ws.Range("J23").Value = 3
$ws.Range("M23").Value = 1.11
$ws.Range("N23").Value = 6.5
